# Append one new data row to the sheet, mirroring the latest Adafruit IO
# sensor reading (same shape/values as the rows already present).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

$timestamp = "2024-09-25T18:06:40Z"
$feedKey   = "temperature"
$value     = "25"
$lat       = "N/A"
$lon       = "N/A"
$elev      = "N/A"

$ws.Cells.Item($newRow, 1).Value = $timestamp
$ws.Cells.Item($newRow, 2).Value = $feedKey

# "25" looks numeric, so a plain assignment would be auto-coerced into a
# Double. Force it to stay text (matching the source data, which stores it
# as a string) the same way Excel does when you type a leading apostrophe,
# then drop back to the workbook's default style so no new number format /
# cell style is introduced.
$ws.Cells.Item($newRow, 3).Value = "'" + $value
$ws.Cells.Item($newRow, 3).Style = "Normal"

$ws.Cells.Item($newRow, 4).Value = $lat
$ws.Cells.Item($newRow, 5).Value = $lon
$ws.Cells.Item($newRow, 6).Value = $elev
